$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.42%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-5.96%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.068"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.67%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07728"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.11%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.341"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.50%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.900"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-7.52%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.185"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.35%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.73%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9200"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.22%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1236"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-8.85%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1868"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.27%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08793"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.21%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03419"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.26%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09700"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.22%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001366"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.36%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006008"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.92%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.569"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.25%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1269"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.05%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.04%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.48%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.02103"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "5,153.47%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04331"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.66%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001211"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.24%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004227"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-11.93%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001350"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.69%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02172"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.54%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04899"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.58%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007657"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.16%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009945"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.99%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1337"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.07%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001994"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.59%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009849"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "5.75%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006552"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.80%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003001"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "3.56%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001301"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-23.08%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"

$ws.Range("D2:E51").ClearFormats()
